# Update DS final project deck:
#  1. Add two new (blank) slides at positions 2 and 3.
#  2. Refresh the cached "datetimeFigureOut" date field text (last-saved
#     date) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# --- 1. Add two new blank slides -------------------------------------------------
# ppLayoutBlank = 12
$p.Slides.Add(2, 12) | Out-Null
$p.Slides.Add(3, 12) | Out-Null

# --- 2. Refresh the cached date placeholder text ---------------------------------
$oldDate = "11/13/2017"
$newDate = "11/26/2017"

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}
